$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Notes placeholder text in C2/C3: "notes" -> "note"
$ws.Range("C3").Value = "{d.records[i+1].note}"
$ws.Range("C2").Value = "{d.records[i].note}"

# Update the selection seen in the saved file (cosmetic, matches target diff)
$ws.Range("E8").Select()
